$wb = $excel.ActiveWorkbook

# Sheet: list
$ws = $wb.Worksheets.Item("list")
$ws.Range("D2").Value = 680
$ws.Range("E2").Value = 644
$ws.Range("F2").Value = 2.94
$ws.Range("H2").Value = 1.91
$ws.Range("J2").Value = 5.29
$ws.Range("L2").Value = 714
$ws.Range("M2").Value = 689
$ws.Range("R2").Value = 3.5
$ws.Range("T2").Value = 706
$ws.Range("U2").Value = 652
$ws.Range("V2").Value = 3.4
$ws.Range("W2").Value = 1.56
$ws.Range("X2").Value = 2.69
$ws.Range("Z2").Value = 7.65
$ws.Range("D3").Value = 680
$ws.Range("E3").Value = 624
$ws.Range("F3").Value = 2.94
$ws.Range("H3").Value = 4.12
$ws.Range("J3").Value = 8.24
$ws.Range("L3").Value = 714
$ws.Range("M3").Value = 681
$ws.Range("N3").Value = 3.5
$ws.Range("R3").Value = 4.62
$ws.Range("T3").Value = 706
$ws.Range("U3").Value = 648
$ws.Range("V3").Value = 1.84
$ws.Range("W3").Value = 0.85
$ws.Range("X3").Value = 5.52
$ws.Range("Z3").Value = 8.22
$ws.Range("D4").Value = 680
$ws.Range("E4").Value = 575
$ws.Range("F4").Value = 3.24
$ws.Range("H4").Value = 11.47
$ws.Range("J4").Value = 15.44
$ws.Range("L4").Value = 714
$ws.Range("M4").Value = 689
$ws.Range("P4").Value = 0.28
$ws.Range("R4").Value = 3.5
$ws.Range("T4").Value = 706
$ws.Range("U4").Value = 598
$ws.Range("V4").Value = 3.26
$ws.Range("X4").Value = 11.76
$ws.Range("Z4").Value = 15.3
$ws.Range("D5").Value = 680
$ws.Range("E5").Value = 505
$ws.Range("F5").Value = 3.24
$ws.Range("G5").Value = 0.59
$ws.Range("H5").Value = 21.91
$ws.Range("J5").Value = 25.74
$ws.Range("L5").Value = 714
$ws.Range("M5").Value = 687
$ws.Range("O5").Value = 1.26
$ws.Range("P5").Value = 0.98
$ws.Range("R5").Value = 3.78
$ws.Range("T5").Value = 706
$ws.Range("U5").Value = 524
$ws.Range("V5").Value = 3.82
$ws.Range("W5").Value = 0.57
$ws.Range("X5").Value = 21.39
$ws.Range("Z5").Value = 25.78
$ws.Range("D6").Value = 680
$ws.Range("E6").Value = 349
$ws.Range("F6").Value = 1.91
$ws.Range("G6").Value = 0.74
$ws.Range("H6").Value = 46.03
$ws.Range("J6").Value = 48.68
$ws.Range("L6").Value = 714
$ws.Range("M6").Value = 669
$ws.Range("N6").Value = 4.06
$ws.Range("P6").Value = 1.54
$ws.Range("R6").Value = 6.3
$ws.Range("T6").Value = 706
$ws.Range("U6").Value = 361
$ws.Range("V6").Value = 2.12
$ws.Range("W6").Value = 0.99
$ws.Range("X6").Value = 45.75
$ws.Range("Z6").Value = 48.87
$ws.Range("D7").Value = 680
$ws.Range("E7").Value = 657
$ws.Range("F7").Value = 2.35
$ws.Range("J7").Value = 3.38
$ws.Range("D8").Value = 680
$ws.Range("E8").Value = 655
$ws.Range("F8").Value = 2.94
$ws.Range("D9").Value = 680
$ws.Range("E9").Value = 647
$ws.Range("H9").Value = 0.29
$ws.Range("J9").Value = 4.85
$ws.Range("D10").Value = 680
$ws.Range("E10").Value = 648
$ws.Range("H10").Value = 0.59
$ws.Range("J10").Value = 4.71
$ws.Range("D11").Value = 680
$ws.Range("E11").Value = 648
$ws.Range("F11").Value = 2.35
$ws.Range("H11").Value = 1.32
$ws.Range("L12").Value = 714
$ws.Range("M12").Value = 676
$ws.Range("N12").Value = 2.8
$ws.Range("P12").Value = 2.1
$ws.Range("R12").Value = 5.32
$ws.Range("L13").Value = 714
$ws.Range("M13").Value = 643
$ws.Range("N13").Value = 3.78
$ws.Range("O13").Value = 1.12
$ws.Range("P13").Value = 5.04
$ws.Range("R13").Value = 9.94
$ws.Range("L14").Value = 714
$ws.Range("M14").Value = 623
$ws.Range("N14").Value = 2.24
$ws.Range("P14").Value = 9.66
$ws.Range("R14").Value = 12.75
$ws.Range("L15").Value = 714
$ws.Range("N15").Value = 3.22
$ws.Range("O15").Value = 0.42
$ws.Range("P15").Value = 16.95
$ws.Range("R15").Value = 20.59
$ws.Range("L16").Value = 714
$ws.Range("N16").Value = 2.38
$ws.Range("O16").Value = 0.42
$ws.Range("P16").Value = 29.27
$ws.Range("R16").Value = 32.07
$ws.Range("L17").Value = 714
$ws.Range("M17").Value = 336
$ws.Range("N17").Value = 1.26
$ws.Range("O17").Value = 0.56
$ws.Range("P17").Value = 51.12
$ws.Range("R17").Value = 52.94
$ws.Range("T18").Value = 706
$ws.Range("U18").Value = 671
$ws.Range("V18").Value = 3.82
$ws.Range("Z18").Value = 4.96
$ws.Range("T19").Value = 706
$ws.Range("U19").Value = 678
$ws.Range("Z19").Value = 3.97
$ws.Range("T20").Value = 706
$ws.Range("U20").Value = 678
$ws.Range("V20").Value = 3.26
$ws.Range("X20").Value = 0.28
$ws.Range("Z20").Value = 3.97
$ws.Range("T21").Value = 706
$ws.Range("U21").Value = 676
$ws.Range("V21").Value = 2.69
$ws.Range("X21").Value = 0.85
$ws.Range("Z21").Value = 4.25
$ws.Range("T22").Value = 706
$ws.Range("U22").Value = 676
$ws.Range("V22").Value = 2.55
$ws.Range("X22").Value = 1.27
$ws.Range("Z22").Value = 4.25

# Sheet: summary_all
$ws = $wb.Worksheets.Item("summary_all")
$ws.Range("E2").Value = 714
$ws.Range("B3").Value = 1714.45
$ws.Range("C3").Value = 522.85
$ws.Range("D3").Value = 1972
$ws.Range("E3").Value = 336
$ws.Range("B4").Value = 2.75
$ws.Range("C4").Value = 0.61
$ws.Range("D4").Value = 2.95
$ws.Range("E4").Value = 1.26
$ws.Range("F4").Value = 3.43
$ws.Range("B5").Value = 0.78
$ws.Range("C5").Value = 0.19
$ws.Range("E5").Value = 0.52
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 17.9
$ws.Range("D6").Value = 2.24
$ws.Range("F6").Value = 51.12
$ws.Range("B8").Value = 15.54
$ws.Range("C8").Value = 17.32
$ws.Range("D8").Value = 6.1
$ws.Range("F8").Value = 52.94

# Sheet: summary_booklet1
$ws = $wb.Worksheets.Item("summary_booklet1")
$ws.Range("B2").Value = 680
$ws.Range("D2").Value = 680
$ws.Range("E2").Value = 680
$ws.Range("F2").Value = 680
$ws.Range("B3").Value = 595.2
$ws.Range("C3").Value = 98.9
$ws.Range("D3").Value = 645.5
$ws.Range("E3").Value = 349
$ws.Range("F3").Value = 657
$ws.Range("B4").Value = 2.87
$ws.Range("C4").Value = 0.52
$ws.Range("D4").Value = 2.94
$ws.Range("E4").Value = 1.91
$ws.Range("B5").Value = 0.84
$ws.Range("C5").Value = 0.23
$ws.Range("E5").Value = 0.44
$ws.Range("B6").Value = 8.76
$ws.Range("C6").Value = 14.85
$ws.Range("D6").Value = 1.61
$ws.Range("F6").Value = 46.03
$ws.Range("B8").Value = 12.47
$ws.Range("C8").Value = 14.54
$ws.Range("D8").Value = 5.07
$ws.Range("E8").Value = 3.38
$ws.Range("F8").Value = 48.68

# Sheet: summary_booklet2
$ws = $wb.Worksheets.Item("summary_booklet2")
$ws.Range("B2").Value = 714
$ws.Range("D2").Value = 714
$ws.Range("E2").Value = 714
$ws.Range("F2").Value = 714
$ws.Range("B3").Value = 613.18
$ws.Range("C3").Value = 111.8
$ws.Range("D3").Value = 669
$ws.Range("E3").Value = 336
$ws.Range("F3").Value = 689
$ws.Range("C4").Value = 0.88
$ws.Range("E4").Value = 1.26
$ws.Range("F4").Value = 4.06
$ws.Range("B5").Value = 0.76
$ws.Range("C5").Value = 0.3
$ws.Range("E5").Value = 0.42
$ws.Range("B6").Value = 10.63
$ws.Range("C6").Value = 16.25
$ws.Range("D6").Value = 2.1
$ws.Range("F6").Value = 51.12
$ws.Range("B8").Value = 14.12
$ws.Range("C8").Value = 15.66
$ws.Range("D8").Value = 6.3
$ws.Range("E8").Value = 3.5
$ws.Range("F8").Value = 52.94

# Sheet: summary_booklet3
$ws = $wb.Worksheets.Item("summary_booklet3")
$ws.Range("B2").Value = 706
$ws.Range("D2").Value = 706
$ws.Range("E2").Value = 706
$ws.Range("F2").Value = 706
$ws.Range("B3").Value = 616.2
$ws.Range("C3").Value = 102.22
$ws.Range("D3").Value = 661.5
$ws.Range("E3").Value = 361
$ws.Range("F3").Value = 678
$ws.Range("B4").Value = 2.97
$ws.Range("C4").Value = 0.67
$ws.Range("D4").Value = 3.12
$ws.Range("E4").Value = 1.84
$ws.Range("F4").Value = 3.82
$ws.Range("C5").Value = 0.39
$ws.Range("D5").Value = 0.78
$ws.Range("F5").Value = 1.56
$ws.Range("B6").Value = 8.95
$ws.Range("C6").Value = 14.64
$ws.Range("D6").Value = 1.98
$ws.Range("F6").Value = 45.75
$ws.Range("B8").Value = 12.72
$ws.Range("C8").Value = 14.48
$ws.Range("D8").Value = 6.3
$ws.Range("E8").Value = 3.97
$ws.Range("F8").Value = 48.87
